$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 4 ("Reference Strain") above the current "Timepoint Type" row,
# pushing Timepoint Type / Cell Location / Value Unit / Scale / Header Format down by one.
$ws.Rows.Item(4).Insert()

# New row 4 content
$ws.Cells.Item(4, 1).Value = "Reference Strain"
$ws.Cells.Item(4, 3).Value = "The Reference Strain (for relative quantification data sets, leave empty for absolute)"

# Match formatting of the surrounding template rows (Insert() only partially
# carries formatting across, so re-apply explicitly from sibling cells).
$ws.Cells.Item(5, 3).Copy() | Out-Null
$ws.Cells.Item(4, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> s="2" like other Description cells

$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Cells.Item(4, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> s="5" like other empty Value cells

$ws.Application.CutCopyMode = $false

# Update the "Value Unit" description text (now on row 7 after the insert).
$ws.Cells.Item(7, 3).Value = "One of mM, uM, Percent, RatioT1, RatioCs, AU, Dimensionless, fmol/ug"

# Insert() carries the old row 8's explicit "ht=18/customHeight" marker down onto
# the new row 9 (Header Format). The finalized template keeps that marker on row 8
# (Scale) instead, so clear it off row 9 and re-apply it to row 8.
$ws.Rows.Item(9).AutoFit() | Out-Null
$ws.Rows.Item(8).RowHeight = 18

# Widen column A slightly to fit the new "Reference Strain" label.
$ws.Columns.Item(1).ColumnWidth = 18.15

# Restore the selection to match the finalized template.
$ws.Range("C18").Select() | Out-Null
